# Updated cryptos list on Wed Mar  1 08:52:53 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row,
# and correct the ordering of the Algorand / InternetComputer(DFINITY) rows
# (rows 37-38) along with their refreshed figures.
#
# NumberFormat is forced to Text ("@") before writing any Price value that
# would otherwise be auto-parsed as a number by Excel, so values such as
# "0.9995" or "1.000" are preserved verbatim as text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.759.11'
$ws.Range('E2').Value = '  +1.96%  '
$ws.Range('D3').Value = '1.653.54'
$ws.Range('E3').Value = '  +1.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9998'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '304.07'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3830'
$ws.Range('E7').Value = '  +2.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3615'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.32'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('E10').Value = '  +2.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08241'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.72'
$ws.Range('E13').Value = '  +2.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.553'
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.419'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '1.653.18'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.78'
$ws.Range('E18').Value = '  +4.33%  '
$ws.Range('E19').Value = '  +0.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.780'
$ws.Range('E20').Value = '  +4.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.72'
$ws.Range('E21').Value = '  +1.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9994'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.66'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('D24').Value = '23.746.92'
$ws.Range('E24').Value = '  +1.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.535'
$ws.Range('E25').Value = '  +3.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.077'
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('E27').Value = '  +1.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.19'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.244'
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '135.01'
$ws.Range('E30').Value = '  +1.97%  '
$ws.Range('D31').Value = '1.835.06'
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.874'
$ws.Range('E32').Value = '  +2.75%  '
$ws.Range('E33').Value = '  +4.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.96'
$ws.Range('E34').Value = '  +10.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.106'
$ws.Range('E35').Value = '  -0.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02841'
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.159'
$ws.Range('E37').Value = '  +3.36%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2521'
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.08841'
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07053'
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.95'
$ws.Range('E41').Value = '  +7.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7084'
$ws.Range('E42').Value = '  +2.02%  '
$ws.Range('E43').Value = '  +0.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.96'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6568'
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.339'
$ws.Range('E46').Value = '  +3.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9992'
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.977'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07994'
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '129.32'
$ws.Range('E50').Value = '  +3.12%  '
$ws.Range('E51').Value = '  +1.59%  '
